$wb = $excel.ActiveWorkbook

# --- Sheet "Valori": rename the first block of POI-type tags (row 1) and
#     extend it with a few more categories (reusing existing tags where
#     they already exist elsewhere in the workbook). ---
$valori = $wb.Worksheets.Item("Valori")

$valori.Range("A1").Value = "Grotta"
$valori.Range("B1").Value = "Ristoro"
$valori.Range("C1").Value = "Accoglienza"
$valori.Range("D1").Value = "Svago"
$valori.Range("E1").Value = "Infopoint"
$valori.Range("F1").Value = "Servizi"
$valori.Range("G1").Value = "Trasporti"
$valori.Range("H1").Value = "Sanità"
$valori.Range("I1").Value = "Segnaletica"
$valori.Range("J1").Value = "Attrazione Naturalistica"
$valori.Range("K1").Value = "Museo"
$valori.Range("L1").Value = "Monumento"

# Reset the remembered selection on that sheet back to A1.
$valori.Range("A1").Select() | Out-Null

# --- Make "QRcode" the active tab (it was "Sito" before). ---
$qrcode = $wb.Worksheets.Item("QRcode")
$qrcode.Activate()
